$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 (Marking): correct the "Wrong" marking count/value
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

# Row 12 (Total): correct the total marks obtained and the displayed "obtained / max" text
$ws.Range("B12").Value = 64
$ws.Range("C12").Value = -2
$ws.Range("E12").Value = "62 / 112"
